$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column P (year 2022) by extending the existing O-column formatting ---

# P3: bottom border cell (same formatting as N3/O3), stays empty
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

# P4: new year header "2022" (same formatting as O4)
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)
$ws.Range("P4").Value = 2022

# P5: new data value for 2022 (same formatting as O5)
$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)
$ws.Range("P5").Value = 2.6

# --- Update the existing 2019-2021 values ---
$ws.Range("M5").Value = 2.6
$ws.Range("N5").Value = 2.4
$ws.Range("O5").Value = 3.3

# --- Move the active selection to P3 ---
[void]$ws.Range("P3").Select()
